{"js": "const replacements = [\n  [\"2025-07-05 Saturday\", \"2025-07-06 Sunday\"],\n  [\"77\u00d779=6083\", \"30\u00d711=330\"],\n  [\"94\u00d793=8742\", \"22\u00d727=594\"],\n  [\"86\u00d731=2666\", \"95\u00d737=3515\"],\n  [\"91\u00d723=2093\", \"65\u00d727=1755\"],\n  [\"31\u00d772=2232\", \"31\u00d797=3007\"],\n  [\"69\u00d723=1587\", \"77\u00d738=2926\"],\n  [\"33\u00d750=1650\", \"63\u00d782=5166\"],\n  [\"64\u00d776=4864\", \"19\u00d782=1558\"],\n  [\"82\u00d720=1640\", \"69\u00d757=3933\"],\n  [\"31\u00d715=465\", \"57\u00d770=3990\"],\n  [\"73\u00d745=3285\", \"17\u00d749=833\"],\n  [\"58\u00d770=4060\", \"71\u00d750=3550\"],\n  [\"23\u00d739=897\", \"19\u00d721=399\"],\n  [\"86\u00d748=4128\", \"65\u00d738=2470\"],\n  [\"56\u00d750=2800\", \"99\u00d773=7227\"],\n  [\"57\u00d711=627\", \"92\u00d723=2116\"],\n  [\"65\u00d763=4095\", \"15\u00d752=780\"],\n  [\"24\u00d797=2328\", \"28\u00d739=1092\"],\n  [\"31\u00d738=1178\", \"90\u00d764=5760\"],\n  [\"29\u00d770=2030\", \"50\u00d775=3750\"],\n  [\"27\u00d797=2619\", \"43\u00d733=1419\"],\n  [\"60\u00d788=5280\", \"42\u00d784=3528\"],\n  [\"81\u00d763=5103\", \"70\u00d728=1960\"],\n  [\"58\u00d767=3886\", \"32\u00d727=864\"],\n  [\"12\u00d789=1068\", \"31\u00d784=2604\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @{Before = \"2025-07-05 Saturday\"; After = \"2025-07-06 Sunday\"},\n    @{Before = \"77\u00d779=6083\"; After = \"30\u00d711=330\"},\n    @{Before = \"94\u00d793=8742\"; After = \"22\u00d727=594\"},\n    @{Before = \"86\u00d731=2666\"; After = \"95\u00d737=3515\"},\n    @{Before = \"91\u00d723=2093\"; After = \"65\u00d727=1755\"},\n    @{Before = \"31\u00d772=2232\"; After = \"31\u00d797=3007\"},\n    @{Before = \"69\u00d723=1587\"; After = \"77\u00d738=2926\"},\n    @{Before = \"33\u00d750=1650\"; After = \"63\u00d782=5166\"},\n    @{Before = \"64\u00d776=4864\"; After = \"19\u00d782=1558\"},\n    @{Before = \"82\u00d720=1640\"; After = \"69\u00d757=3933\"},\n    @{Before = \"31\u00d715=465\"; After = \"57\u00d770=3990\"},\n    @{Before = \"73\u00d745=3285\"; After = \"17\u00d749=833\"},\n    @{Before = \"58\u00d770=4060\"; After = \"71\u00d750=3550\"},\n    @{Before = \"23\u00d739=897\"; After = \"19\u00d721=399\"},\n    @{Before = \"86\u00d748=4128\"; After = \"65\u00d738=2470\"},\n    @{Before = \"56\u00d750=2800\"; After = \"99\u00d773=7227\"},\n    @{Before = \"57\u00d711=627\"; After = \"92\u00d723=2116\"},\n    @{Before = \"65\u00d763=4095\"; After = \"15\u00d752=780\"},\n    @{Before = \"24\u00d797=2328\"; After = \"28\u00d739=1092\"},\n    @{Before = \"31\u00d738=1178\"; After = \"90\u00d764=5760\"},\n    @{Before = \"29\u00d770=2030\"; After = \"50\u00d775=3750\"},\n    @{Before = \"27\u00d797=2619\"; After = \"43\u00d733=1419\"},\n    @{Before = \"60\u00d788=5280\"; After = \"42\u00d784=3528\"},\n    @{Before = \"81\u00d763=5103\"; After = \"70\u00d728=1960\"},\n    @{Before = \"58\u00d767=3886\"; After = \"32\u00d727=864\"},\n    @{Before = \"12\u00d789=1068\"; After = \"31\u00d784=2604\"},\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Before, $false, $true, $false, $false, $false, $true, 1, $false, $r.After, 2)\n}\n"}
